$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310".
# Columns A:J used the "_old" suffix, column K is the untouched "diff" header,
# and columns L:U used the "_new" suffix.
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = [string]$cell.Value2 -replace "_old$", "_FV2304"
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = [string]$cell.Value2 -replace "_new$", "_FV2310"
}

# Turn the used range into an Excel Table ("ListObject") with AutoFilter.
$rng = $ws.Range("A1:U56")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1) and keep the selection on the
# scrollable pane, matching the sheetView/pane added to the worksheet.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
